$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I to make room for "Image Filename"
$ws.Columns("I:I").Insert()

# Header row: set the new column header
$ws.Range("I1").Value = "Image Filename"

# --- Row 2 (Name0) updates ---
$ws.Range("B2").Value = "Username0"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "23"

$ws.Range("H2").Value = "Battambang"
$ws.Range("I2").Value = "Name0_Username0_20250422_225605.jpg"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "81991"

$ws.Range("K2").Value = "22/04/2025"
$ws.Range("L2").Value = "22:56:05"
$ws.Range("M2").Value = "B821"

# --- Row 3 (Name1) updates ---
$ws.Range("B3").Value = "UserName1"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "24"

$ws.Range("D3").Value = "Male"
$ws.Range("I3").Value = "Name1_UserName1_20250422_225703.png"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "43575"

$ws.Range("K3").Value = "22/04/2025"
$ws.Range("L3").Value = "22:57:03"
$ws.Range("M3").Value = "B944"

# --- Remove old rows 4-6 (students Name3, Name4, Name5 are no longer present) ---
$ws.Rows("4:6").Delete()
